$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3380.4
$ws.Range("I40").Value = 1898
$ws.Range("J40").Value = 3751
$ws.Range("K40").Value = 1898
$ws.Range("L40").Value = 3751
$ws.Range("M40").Value = -1723
$ws.Range("N40").Value = -4101
$ws.Range("H46").Value = 1076
$ws.Range("I46").Value = 935
$ws.Range("K46").Value = 2805
$ws.Range("M46").Value = -2686
$ws.Range("H54").Value = 19000
$ws.Range("I54").Value = 21000
$ws.Range("K54").Value = 21000
$ws.Range("M54").Value = -20514
$ws.Range("H60").Value = 1076
$ws.Range("I60").Value = 935
$ws.Range("K60").Value = 2805
$ws.Range("M60").Value = -2321
$ws.Range("H76").Value = 2530204
$ws.Range("I76").Value = 4278262.5
$ws.Range("K76").Value = 4278262.5
$ws.Range("M76").Value = -4277947.5
$ws.Range("H79").Value = 2530204
$ws.Range("I79").Value = 4278262.5
$ws.Range("K79").Value = 4278262.5
$ws.Range("M79").Value = -4277170.5
$ws.Range("H86").Value = 5294.467
$ws.Range("I86").Value = 5610.5557
$ws.Range("K86").Value = 5610.5557
$ws.Range("M86").Value = -4487.5557
$ws.Range("H89").Value = 5294.467
$ws.Range("I89").Value = 5610.5557
$ws.Range("K89").Value = 28052.7785
$ws.Range("M89").Value = -22436.7785
$ws.Range("H132").Value = 26318852
$ws.Range("I132").Value = 28574608
$ws.Range("K132").Value = 85723824
$ws.Range("M132").Value = -85721294
$ws.Range("H137").Value = 3122.3076
$ws.Range("I137").Value = 2116.3157
$ws.Range("K137").Value = 6348.9471
$ws.Range("M137").Value = -3798.9471
$ws.Range("H138").Value = 1929.3131
$ws.Range("J138").Value = 2807.7407
$ws.Range("L138").Value = 8423.222099999999
$ws.Range("N138").Value = -18703.2221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 672.08
$ws.Range("I32").Value = 600.0833
$ws.Range("K32").Value = 600.0833
$ws.Range("M32").Value = -313.0833
$ws.Range("H102").Value = 2691153.2
$ws.Range("I102").Value = 3336036
$ws.Range("K102").Value = 3336036
$ws.Range("M102").Value = -3334414
$ws.Range("H132").Value = 2484.6667
$ws.Range("J132").Value = 3503.2
$ws.Range("L132").Value = 10509.6
$ws.Range("N132").Value = -15569.6
$ws.Range("H137").Value = 80259.664
$ws.Range("J137").Value = 70389.5
$ws.Range("L137").Value = 70389.5
$ws.Range("N137").Value = -80589.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2827.5386
$ws.Range("I20").Value = 2481
$ws.Range("J20").Value = 3382
$ws.Range("K20").Value = 2481
$ws.Range("L20").Value = 3382
$ws.Range("M20").Value = -2234
$ws.Range("N20").Value = -3876
$ws.Range("H107").Value = 7937973
$ws.Range("I107").Value = 7937973
$ws.Range("K107").Value = 7937973
$ws.Range("M107").Value = -7936053
$ws.Range("H120").Value = 99999
$ws.Range("J120").Value = 99999
$ws.Range("L120").Value = 99999
$ws.Range("N120").Value = -109675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 12000.333
$ws.Range("I86").Value = 10251.75
$ws.Range("J86").Value = 13399.2
$ws.Range("K86").Value = 10251.75
$ws.Range("L86").Value = 13399.2
$ws.Range("M86").Value = -9128.75
$ws.Range("N86").Value = -15645.2
$ws.Range("H89").Value = 12000.333
$ws.Range("I89").Value = 10251.75
$ws.Range("J89").Value = 13399.2
$ws.Range("K89").Value = 51258.75
$ws.Range("L89").Value = 66996
$ws.Range("M89").Value = -45642.75
$ws.Range("N89").Value = -78228
$ws.Range("H99").Value = 3668.8572
$ws.Range("I99").Value = 2736.75
$ws.Range("K99").Value = 2736.75
$ws.Range("M99").Value = -1238.75
$ws.Range("H126").Value = 3668.8572
$ws.Range("I126").Value = 2736.75
$ws.Range("K126").Value = 8210.25
$ws.Range("M126").Value = -5740.25
$ws.Range("H132").Value = 45390.39
$ws.Range("I132").Value = 57110.332
$ws.Range("J132").Value = 3198.6
$ws.Range("K132").Value = 171330.996
$ws.Range("L132").Value = 9595.799999999999
$ws.Range("M132").Value = -168800.996
$ws.Range("N132").Value = -14655.8
$ws.Range("H141").Value = 34320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1075.3334
$ws.Range("J52").Value = 1075.3334
$ws.Range("L52").Value = 3226.0002
$ws.Range("N52").Value = -3758.0002
$ws.Range("H87").Value = 12991.429
$ws.Range("I87").Value = 9388
$ws.Range("K87").Value = 28164
$ws.Range("M87").Value = -26916
$ws.Range("H90").Value = 12991.429
$ws.Range("I90").Value = 9388
$ws.Range("K90").Value = 84492
$ws.Range("M90").Value = -78252
$ws.Range("H124").Value = 1480
$ws.Range("I124").Value = 1150
$ws.Range("K124").Value = 3450
$ws.Range("M124").Value = 1460
$ws.Range("H132").Value = 2776.7896
$ws.Range("I132").Value = 2454.2856
$ws.Range("J132").Value = 2964.9167
$ws.Range("K132").Value = 22088.5704
$ws.Range("L132").Value = 26684.2503
$ws.Range("M132").Value = -19558.5704
$ws.Range("N132").Value = -31744.2503
$ws.Range("H137").Value = 2769.375
$ws.Range("J137").Value = 3720.1428
$ws.Range("L137").Value = 11160.4284
$ws.Range("N137").Value = -21360.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7578880
$ws.Range("I126").Value = 3790983.5
$ws.Range("J126").Value = 16669831
$ws.Range("K126").Value = 11372950.5
$ws.Range("L126").Value = 50009493
$ws.Range("M126").Value = -11370480.5
$ws.Range("N126").Value = -50014433

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5028.852
$ws.Range("I40").Value = 4204.7896
$ws.Range("K40").Value = 4204.7896
$ws.Range("M40").Value = -4068.7896
$ws.Range("H132").Value = 4908.3433
$ws.Range("I132").Value = 4800.5
$ws.Range("K132").Value = 14401.5
$ws.Range("M132").Value = -11871.5
$ws.Range("H136").Value = 47919.89
$ws.Range("I136").Value = 72232.86
$ws.Range("K136").Value = 216698.58
$ws.Range("M136").Value = -214148.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 8574
$ws.Range("J31").Value = 9169.666999999999
$ws.Range("L31").Value = 9169.666999999999
$ws.Range("N31").Value = -9865.666999999999
$ws.Range("H116").Value = 75000
$ws.Range("J116").Value = 75000
$ws.Range("L116").Value = 75000
$ws.Range("N116").Value = -84178
